# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet (the active sheet), a new blank column is
# inserted immediately before the existing "Late" column. That pushes the old
# N/O/P columns (Late / blank / Outstanding) one slot to the right -> O/P/Q,
# with a fresh blank column taking over the old "N" slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column N ("Late"), shifting
# Late / blank / Outstanding from N/O/P to O/P/Q.
$ws.Columns("N").Insert() | Out-Null

# Leave the selection where it ended up after making the edit.
$ws.Range("L17").Select() | Out-Null
